$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "26.259.04"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.65%  "

$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("E4").Value = "  +0.68%  "

$ws.Range("E5").Value = "  +0.16%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.5327"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.94%  "

$ws.Range("E7").Value = "  +0.62%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2633"
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06346"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.68%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "20.50"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.78%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07840"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.37%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "4.540"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.64%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.641.46"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.54%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "1.888.36"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.57%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.5516"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.26%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0₅8178"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.88%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "65.60"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.70%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "26.247.63"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.57%  "

$ws.Range("E19").Value = "  +0.66%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "4.651"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.44%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "191.94"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.80%  "

$ws.Range("E22").Value = "  +0.86%  "

$ws.Range("E23").Value = "  +1.26%  "

$ws.Range("E24").Value = "  +0.64%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "144.72"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +3.37%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.1228"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.86%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "7.226"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.10%  "

$ws.Range("E28").Value = "  -0.37%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.470"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.49%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.05783"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.88%  "

$ws.Range("E31").Value = "  +0.11%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.572"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.09%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.286"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.62%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.611"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +4.14%  "

$ws.Range("E35").Value = "  +2.16%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.9568"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.58%  "

$ws.Range("E37").Value = "  +0.66%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.5786"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +2.55%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01603"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.13%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "5.844"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.05%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.8517"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.88%  "

$ws.Range("E42").Value = "  +0.61%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "104.71"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +3.95%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.045.77"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +3.90%  "

$ws.Range("E45").Value = "  +0.47%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "57.06"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.50%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0₈106"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.17%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.013"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.66%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.4370"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.91%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "7.991"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.35%  "

$ws.Range("E51").Value = "  +0.20%  "
